# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new values look like plain numbers (e.g. "4.70", "0.998") ---
# Force them to remain TEXT (matching the source inlineStr cells) by switching
# the cell to a text number-format before assigning, then restoring the default
# "Normal" style so no stray formatting is left behind.
$textForcedCells = @('D5', 'D6', 'D8', 'D10', 'D13', 'D15', 'D16', 'D18', 'D19', 'D21', 'D22', 'D23', 'D25', 'D27', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D40', 'D42', 'D46', 'D49', 'D50')
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D5').Value = '227.63'
$ws.Range('D6').Value = '0.611'
$ws.Range('D8').Value = '42.31'
$ws.Range('D10').Value = '0.0691'
$ws.Range('D13').Value = '11.51'
$ws.Range('D15').Value = '4.75'
$ws.Range('D16').Value = '0.659'
$ws.Range('D18').Value = '69.52'
$ws.Range('D19').Value = '245.32'
$ws.Range('D21').Value = '12.05'
$ws.Range('D22').Value = '4.70'
$ws.Range('D23').Value = '0.998'
$ws.Range('D25').Value = '171.52'
$ws.Range('D27').Value = '17.85'
$ws.Range('D31').Value = '3.95'
$ws.Range('D32').Value = '4.03'
$ws.Range('D33').Value = '0.0536'
$ws.Range('D34').Value = '1.90'
$ws.Range('D35').Value = '0.671'
$ws.Range('D36').Value = '90.16'
$ws.Range('D40').Value = '2.42'
$ws.Range('D42').Value = '14.96'
$ws.Range('D46').Value = '0.0517'
$ws.Range('D49').Value = '103.97'
$ws.Range('D50').Value = '0.998'

foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Remaining cells: plain text / URLs / percentages / multi-dot numbers ---
# These are not number-like so the COM layer stores them as text natively.
$ws.Range('D2').Value = '35.188.27'
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '1.844.54'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('E6').Value = '  +1.88%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +13.66%  '
$ws.Range('E9').Value = '  +3.99%  '
$ws.Range('E10').Value = '  +1.04%  '
$ws.Range('E11').Value = '  +3.39%  '
$ws.Range('D12').Value = '2.113.16'
$ws.Range('E12').Value = '  +1.90%  '
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('D14').Value = '1.843.88'
$ws.Range('E14').Value = '  +1.46%  '
$ws.Range('E15').Value = '  +7.06%  '
$ws.Range('E16').Value = '  +3.84%  '
$ws.Range('D17').Value = '35.180.87'
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('E21').Value = '  +6.98%  '
$ws.Range('E22').Value = '  +14.07%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E24').Value = '  -1.26%  '
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E27').Value = '  +3.32%  '
$ws.Range('E28').Value = '  +0.65%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +8.47%  '
$ws.Range('E31').Value = '  +3.22%  '
$ws.Range('E32').Value = '  +2.26%  '
$ws.Range('E33').Value = '  +3.46%  '
$ws.Range('E34').Value = '  +4.19%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E35').Value = '  +2.48%  '
$ws.Range('B36').Value = 'Aave'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E36').Value = '  +11.06%  '
$ws.Range('E37').Value = '  +1.97%  '
$ws.Range('D38').Value = '1.340.59'
$ws.Range('E38').Value = '  -1.87%  '
$ws.Range('E39').Value = '  +8.98%  '
$ws.Range('E40').Value = '  +1.73%  '
$ws.Range('E41').Value = '  +2.89%  '
$ws.Range('E42').Value = '  +8.75%  '
$ws.Range('E43').Value = '  +6.70%  '
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('E46').Value = '  +3.25%  '
$ws.Range('E47').Value = '  +4.45%  '
$ws.Range('D48').Value = '2.014.43'
$ws.Range('E48').Value = '  +2.03%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0121'
$ws.Range('E51').Value = '  +0.40%  '
